# Update the "Enterprises density (per 1000 people)" row (row 11) and the
# "Enterprises (% of total)" row (row 13) on the Summary sheet with more
# precise figures. The source values are numeric-looking strings that must
# stay stored as TEXT (matching the original file, where they live in
# sharedStrings as <t> entries) rather than being auto-converted to numbers
# by Excel's type inference. A leading apostrophe forces text entry; the
# Style re-assignment afterwards restores the cell's original (default)
# formatting so only the content changes, not the number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $cell = $ws.Range($rangeAddress)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Enterprises density (per 1000 people): Micro / SMEs / MSMEs
Set-TextValue "B11" "22.67"
Set-TextValue "C11" "13.84"
Set-TextValue "D11" "36.51"

# Enterprises (% of total): Micro / SMEs / MSMEs
Set-TextValue "B13" "61.62"
Set-TextValue "C13" "37.61"
Set-TextValue "D13" "99.23"
